# Auto-generated Excel COM-interop script applying the Asura_Profits market-data refresh diff.
# For each touched cell: set the new value, or clear the cell entirely when the diff removes it.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2350
$ws.Range("H70").Value = 56617
$ws.Range("I70").Value = 334433.34
$ws.Range("J70").Value = 1053.7333
$ws.Range("K70").Value = 1003300.02
$ws.Range("L70").Value = 3161.199900000001
$ws.Range("M70").Value = -1003030.02
$ws.Range("N70").Value = -3701.199900000001
$ws.Range("H73").Value = 56617
$ws.Range("I73").Value = 334433.34
$ws.Range("J73").Value = 1053.7333
$ws.Range("K73").Value = 1003300.02
$ws.Range("L73").Value = 3161.199900000001
$ws.Range("M73").Value = -1002364.02
$ws.Range("N73").Value = -5033.199900000001
$ws.Range("H129").Value = 939.86957
$ws.Range("I129").Value = 498.66666
$ws.Range("J129").Value = 1047.1892
$ws.Range("K129").Value = 1495.99998
$ws.Range("L129").Value = 3141.5676
$ws.Range("M129").Value = 3504.00002
$ws.Range("N129").Value = -13141.5676

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8690.625
$ws.Range("I32").Value = 8910.834999999999
$ws.Range("K32").Value = 8910.834999999999
$ws.Range("M32").Value = -8623.834999999999
$ws.Range("H108").Value = 39700
$ws.Range("J108").Value = 39700
$ws.Range("L108").Value = 39700
$ws.Range("N108").Value = -47380

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 10305.6
$ws.Range("I25").Value = 10305.6
$ws.Range("K25").Value = 10305.6
$ws.Range("M25").Value = -10070.6
$ws.Range("H29").Value = 887.5
$ws.Range("I29").Value = 887.5
$ws.Range("K29").Value = 887.5
$ws.Range("M29").Value = -598.5
$ws.Range("H36").Value = 667.63635
$ws.Range("I36").Value = 667.63635
$ws.Range("K36").Value = 667.63635
$ws.Range("M36").Value = -133.63635
$ws.Range("H107").Value = 500000
$ws.Range("I107").Value = 500000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 500000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -498080
$ws.Range("N107").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16131317
$ws.Range("J31").Value = 4459.467
$ws.Range("L31").Value = 4459.467
$ws.Range("N31").Value = -5049.467
$ws.Range("H34").Value = 16131317
$ws.Range("J34").Value = 4459.467
$ws.Range("L34").Value = 4459.467
$ws.Range("N34").Value = -4863.467
$ws.Range("H107").Value = 692.9091
$ws.Range("I107").Value = 684.4
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 684.4
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = 1235.6
$ws.Range("N107").Value = -4540

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1853.9
$ws.Range("I117").Value = 615
$ws.Range("J117").Value = 2163.625
$ws.Range("K117").Value = 1845
$ws.Range("L117").Value = 6490.875
$ws.Range("M117").Value = 1597
$ws.Range("N117").Value = -13374.875
$ws.Range("H123").Value = 4023.6365
$ws.Range("I123").Value = 1065
$ws.Range("J123").Value = 5714.2856
$ws.Range("K123").Value = 3195
$ws.Range("L123").Value = 17142.8568
$ws.Range("M123").Value = -745
$ws.Range("N123").Value = -22042.8568
$ws.Range("H131").Value = 850.33
$ws.Range("I131").Value = 455.125
$ws.Range("J131").Value = 884.6957
$ws.Range("K131").Value = 1365.375
$ws.Range("L131").Value = 2654.0871
$ws.Range("M131").Value = 3674.625
$ws.Range("N131").Value = -12734.0871
$ws.Range("H133").Value = 3371.6365
$ws.Range("I133").Value = 1326.8572
$ws.Range("J133").Value = 6950
$ws.Range("K133").Value = 3980.5716
$ws.Range("L133").Value = 20850
$ws.Range("M133").Value = 1079.4284
$ws.Range("N133").Value = -30970
$ws.Range("H134").Value = 6778.151
$ws.Range("I134").Value = 6963.45
$ws.Range("K134").Value = 20890.35
$ws.Range("M134").Value = -15820.35

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 557.1429000000001
$ws.Range("I107").Value = 590
$ws.Range("J107").Value = 475
$ws.Range("K107").Value = 590
$ws.Range("L107").Value = 475
$ws.Range("M107").Value = 1330
$ws.Range("N107").Value = -4315
$ws.Range("H113").Value = 1060.5238
$ws.Range("I113").Value = 965.4
$ws.Range("K113").Value = 965.4
$ws.Range("M113").Value = 1204.6

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4254.154
$ws.Range("I7").Value = 3614.8572
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 3614.8572
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -3502.8572
$ws.Range("N7").Value = -5224
$ws.Range("H22").Value = 2410
$ws.Range("I22").Value = 2762.5
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 2762.5
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -2467.5
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 2410
$ws.Range("I27").Value = 2762.5
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 2762.5
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -2655.5
$ws.Range("N27").Value = -1214
$ws.Range("H40").Value = 4658
$ws.Range("I40").Value = 5197.5
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 5197.5
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -5061.5
$ws.Range("N40").Value = -2772
$ws.Range("H126").Value = 4254.154
$ws.Range("I126").Value = 3614.8572
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 10844.5716
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -8374.571599999999
$ws.Range("N126").Value = -19940

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 12475
$ws.Range("I33").Value = 5000
$ws.Range("J33").Value = 14966.667
$ws.Range("K33").Value = 5000
$ws.Range("L33").Value = 14966.667
$ws.Range("M33").Value = -4750
$ws.Range("N33").Value = -15466.667
$ws.Range("H36").Value = 12475
$ws.Range("I36").Value = 5000
$ws.Range("J36").Value = 14966.667
$ws.Range("K36").Value = 5000
$ws.Range("L36").Value = 14966.667
$ws.Range("M36").Value = -4750
$ws.Range("N36").Value = -15466.667
$ws.Range("H107").Value = 486.6842
$ws.Range("I107").Value = 357.75
$ws.Range("K107").Value = 1073.25
$ws.Range("M107").Value = 846.75
